$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value (preserve text formatting,
# matching the original cells which are stored as text/inline strings).
$updates = @{
    "D2" = "304.80"
    "E2" = "-0.36%"
    "D3" = "35.36"
    "E3" = "-2.82%"
    "D4" = "5.108"
    "E4" = "0.95%"
    "D5" = "0.08008"
    "E5" = "1.00%"
    "D6" = "1.948"
    "E6" = "-11.54%"
    "D7" = "7.865"
    "E7" = "-1.89%"
    "D8" = "2.896"
    "E8" = "10.05%"
    "D9" = "0.9238"
    "E9" = "-0.73%"
    "D10" = "0.1091"
    "E10" = "10.70%"
    "D11" = "0.1893"
    "E11" = "0.84%"
    "D12" = "0.09429"
    "E12" = "3.52%"
    "D13" = "0.03703"
    "E13" = "0.19%"
    "E14" = "-0.19%"
    "D15" = "0.001432"
    "E15" = "-0.97%"
    "D16" = "0.005799"
    "E16" = "2.57%"
    "D17" = "3.458"
    "E17" = "-0.22%"
    "D18" = "4.144"
    "E18" = "-0.72%"
    "D19" = "0.3421"
    "E19" = "1.49%"
    "D20" = "0.1321"
    "E20" = "-1.75%"
    "D21" = "5.101"
    "E21" = "0.06%"
    "E22" = "0.32%"
    "D23" = "0.04532"
    "E23" = "-0.53%"
    "E24" = "-0.63%"
    "D25" = "0.004657"
    "E25" = "-2.70%"
    "D26" = "0.0001252"
    "E26" = "-3.71%"
    "D27" = "0.0004488"
    "E27" = "-5.29%"
    "D39" = "0.01875"
    "E39" = "-5.62%"
    "D40" = "0.04748"
    "E40" = "-3.67%"
    "D41" = "0.007544"
    "E41" = "-3.23%"
    "D42" = "0.009689"
    "E42" = "24.09%"
    "D43" = "0.1336"
    "E43" = "-4.40%"
    "D44" = "0.002119"
    "E44" = "0.29%"
    "D45" = "0.01130"
    "E45" = "0.69%"
    "D46" = "0.00006492"
    "E46" = "4.34%"
    "D47" = "0.00000000752"
    "E47" = "0.19%"
    "D48" = "64.38"
    "E48" = "23.27%"
    "D49" = "0.001312"
    "E49" = "-27.17%"
    "D50" = "0.00002105"
    "E50" = "0.19%"
    "D51" = "0.0002005"
    "E51" = "0.19%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text format so Excel does not reinterpret the numeric-looking
    # string (e.g. "304.80") or percent string (e.g. "-0.36%") as a number.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Drop the temporary text-format style so the cell keeps its original
    # (default/general) style, matching the source workbook exactly.
    $cell.ClearFormats()
}
